$wb = $excel.ActiveWorkbook

# Template sheet: "Spain" has the same layout (cols/rows/styles) the new
# Russia/Finland/Hungary sheets are based on (wide "Market" column B,
# "ht=28.8" rows 3-5, style 9 on the code cell B4).
$template = $wb.Worksheets.Item("Spain")

function Add-MarketSheet($name, $code, $market) {
    $lastIndex = $wb.Worksheets.Count
    $afterSheet = $wb.Worksheets.Item($lastIndex)
    $template.Copy($null, $afterSheet)

    $ws = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws.Name = $name

    # Set the code cell before the market-name cell so the shared-string
    # table receives the two new strings in (code, market) order.
    $ws.Range("B4").Value = $code
    $ws.Range("B2").Value = $market

    # Match the narrower column widths used on the new sheets.
    $ws.Columns.Item(2).ColumnWidth = 16
    $ws.Columns.Item(4).ColumnWidth = 7.67

    # New sheets select the whole data block instead of a single cell.
    $ws.Range("A1:D11").Select()

    return $ws
}

Add-MarketSheet "Russia" "NGC-2929/T2925" "Russia Market" | Out-Null
Add-MarketSheet "Finland" "NGC-3130/T2943" "Finland Market" | Out-Null
$hungary = Add-MarketSheet "Hungary" "NGC-3104/T2992" "Hungary Market"

# Hungary is the new last sheet and becomes the active/selected tab.
$hungary.Activate()
